$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")
$ws.Activate()

# Main input change: cost inflation rate assumption drops from 6% to 4%,
# which ripples through the whole projection via the dependent formulas.
$ws.Range("Y30").Value = 0.04

# Formatting tweaks that came with the same edit: the three assumption
# inputs move from 0% to 0.00% precision, and the per-share NPV figure
# switches from the generic "3" style to a 2-decimal number format.
$ws.Range("Y30:Y32").NumberFormat = "0.00%"
$ws.Range("Y34").NumberFormat = "#,##0.00"

# Leave the selection where the author left it.
$ws.Range("Y32").Select()
